# Updates the cryptos price/volume table to the latest scrape.
# D-column values that look like plain numbers get a leading apostrophe
# (Excel's standard "force text" quote-prefix) so they keep being stored
# as text, exactly like the other (non-numeric-looking) price strings that
# already use a '.'-per-thousands format (e.g. "27.920.20").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.919.25'
$ws.Range("D3").Value = '1.642.92'
$ws.Range("E3").Value = '  +1.21%  '
$ws.Range("D5").Value = '''213.45'
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '''23.57'
$ws.Range("E8").Value = '  +1.52%  '
$ws.Range("E9").Value = '  +0.62%  '
$ws.Range("E10").Value = '  +0.81%  '
$ws.Range("D11").Value = '''0.0873'
$ws.Range("E11").Value = '  -1.99%  '
$ws.Range("D12").Value = '1.876.04'
$ws.Range("E12").Value = '  +1.23%  '
$ws.Range("D13").Value = '1.639.33'
$ws.Range("E13").Value = '  +0.95%  '
$ws.Range("D14").Value = '''0.573'
$ws.Range("E14").Value = '  +4.33%  '
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").Value = '''65.78'
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("D17").Value = '27.901.91'
$ws.Range("D18").Value = '''230.40'
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("E19").Value = '  +0.89%  '
$ws.Range("D20").Value = '''7.62'
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D22").Value = '''10.92'
$ws.Range("E22").Value = '  +4.90%  '
$ws.Range("E23").Value = '  +1.51%  '
$ws.Range("D24").Value = '''2.14'
$ws.Range("E24").Value = '  +2.72%  '
$ws.Range("D25").Value = '''152.23'
$ws.Range("E25").Value = '  +1.71%  '
$ws.Range("D26").Value = '''6.91'
$ws.Range("E27").Value = '  +0.84%  '
$ws.Range("D28").Value = '''15.71'
$ws.Range("E28").Value = '  +1.06%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  +0.92%  '
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("E32").Value = '  +2.04%  '
$ws.Range("D33").Value = '1.427.87'
$ws.Range("D34").Value = '''3.09'
$ws.Range("E34").Value = '  +0.82%  '
$ws.Range("E35").Value = '  +1.55%  '
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").Value = '''0.889'
$ws.Range("E37").Value = '  +1.78%  '
$ws.Range("E38").Value = '  +0.64%  '
$ws.Range("E39").Value = '  +0.54%  '
$ws.Range("D40").Value = '''0.921'
$ws.Range("E40").Value = '  -2.38%  '
$ws.Range("E41").Value = '  +2.29%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '''68.61'
$ws.Range("E42").Value = '  +1.58%  '
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = '''1.00'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("E44").Value = '  +0.41%  '
$ws.Range("D45").Value = '''5.44'
$ws.Range("E45").Value = '  +2.96%  '
$ws.Range("E46").Value = '  +2.91%  '
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("D48").Value = '1.784.50'
$ws.Range("E48").Value = '  +1.20%  '
$ws.Range("D49").Value = '''89.08'
$ws.Range("E49").Value = '  +1.89%  '
$ws.Range("E50").Value = '  +0.45%  '
$ws.Range("E51").Value = '  +0.67%  '

Write-Output "Applied cryptos update"
